$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$oldLine1 = "✅ 1000 Bs = 4.63 = 18316.85 pesos"
$newLine1 = "✅ 1000 Bs = 4.74 = 18736.08 pesos"
$oldLine2 = "✅ 18316.85 pesos = 4.59 = 940.82 Bs"
$newLine2 = "✅ 18736.08 pesos = 4.7 = 931.16 Bs"

$text = $ws1.Range("A1").Value()
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$ws1.Range("A1").Value = $text

# --- Sheet "tasas": update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 210.9
$ws2.Range("O10").Value = 3951.44
$ws2.Range("N12").Value = 3984
$ws2.Range("O12").Value = 198
